$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number (45171 -> 2023-09-02).
# Update every data row (2..205) to the new date serial number 45172 (2023-09-03).
$ws.Range("C2:C205").Value = 45172
